$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Sheet ALC, row 10
$ws_ALC.Range("H10").Value = 31333.334
$ws_ALC.Range("J10").Value = 31333.334
$ws_ALC.Range("L10").Value = 31333.334
$ws_ALC.Range("N10").Value = -31919.334

# Sheet ALC, row 129
$ws_ALC.Range("H129").Value = 926.76
$ws_ALC.Range("I129").Value = 351.83334
$ws_ALC.Range("J129").Value = 1108.3158
$ws_ALC.Range("K129").Value = 1055.50002
$ws_ALC.Range("L129").Value = 3324.9474
$ws_ALC.Range("M129").Value = 3944.49998
$ws_ALC.Range("N129").Value = -13324.9474

# Sheet ALC, row 136
$ws_ALC.Range("H136").Value = 0
$ws_ALC.Range("J136").Value = 0
$ws_ALC.Range("L136").Value = 0
$ws_ALC.Range("N136").ClearContents()

# Sheet ALC, row 138
$ws_ALC.Range("H138").Value = 4351606.5
$ws_ALC.Range("I138").Value = 2209.5334
$ws_ALC.Range("J138").Value = 6456153
$ws_ALC.Range("K138").Value = 6628.600199999999
$ws_ALC.Range("L138").Value = 19368459
$ws_ALC.Range("M138").Value = -1488.600199999999
$ws_ALC.Range("N138").Value = -19378739

# Sheet ALC, row 140
$ws_ALC.Range("H140").Value = 0
$ws_ALC.Range("J140").Value = 0
$ws_ALC.Range("L140").Value = 0
$ws_ALC.Range("N140").ClearContents()

$ws_ARM = $wb.Worksheets.Item("ARM")
# Sheet ARM, row 9
$ws_ARM.Range("H9").Value = 39003
$ws_ARM.Range("J9").Value = 17998
$ws_ARM.Range("L9").Value = 17998
$ws_ARM.Range("N9").Value = -18338

# Sheet ARM, row 20
$ws_ARM.Range("H20").Value = 39003
$ws_ARM.Range("J20").Value = 17998
$ws_ARM.Range("L20").Value = 17998
$ws_ARM.Range("N20").Value = -18538

# Sheet ARM, row 42
$ws_ARM.Range("H42").Value = 14960
$ws_ARM.Range("J42").Value = 14960
$ws_ARM.Range("L42").Value = 14960
$ws_ARM.Range("N42").Value = -15932

# Sheet ARM, row 61
$ws_ARM.Range("H61").Value = 2814
$ws_ARM.Range("I61").Value = 1019.1429
$ws_ARM.Range("J61").Value = 4489.2
$ws_ARM.Range("K61").Value = 1019.1429
$ws_ARM.Range("L61").Value = 4489.2
$ws_ARM.Range("M61").Value = -807.1429000000001
$ws_ARM.Range("N61").Value = -4913.2

# Sheet ARM, row 136
$ws_ARM.Range("H136").Value = 2814
$ws_ARM.Range("I136").Value = 1019.1429
$ws_ARM.Range("J136").Value = 4489.2
$ws_ARM.Range("K136").Value = 3057.4287
$ws_ARM.Range("L136").Value = 13467.6
$ws_ARM.Range("M136").Value = -507.4287000000004
$ws_ARM.Range("N136").Value = -18567.6

$ws_BSM = $wb.Worksheets.Item("BSM")
# Sheet BSM, row 25
$ws_BSM.Range("H25").Value = 3504.2727
$ws_BSM.Range("I25").Value = 1093.375
$ws_BSM.Range("J25").Value = 9933.333000000001
$ws_BSM.Range("K25").Value = 1093.375
$ws_BSM.Range("L25").Value = 9933.333000000001
$ws_BSM.Range("M25").Value = -858.375
$ws_BSM.Range("N25").Value = -10403.333

# Sheet BSM, row 107
$ws_BSM.Range("H107").Value = 1878.36
$ws_BSM.Range("I107").Value = 1923.8695
$ws_BSM.Range("J107").Value = 1355
$ws_BSM.Range("K107").Value = 1923.8695
$ws_BSM.Range("L107").Value = 1355
$ws_BSM.Range("M107").Value = -3.869500000000016
$ws_BSM.Range("N107").Value = -5195

$ws_CRP = $wb.Worksheets.Item("CRP")
# Sheet CRP, row 31
$ws_CRP.Range("H31").Value = 7247514.5
$ws_CRP.Range("I31").Value = 699.2069
$ws_CRP.Range("J31").Value = 12501455
$ws_CRP.Range("K31").Value = 699.2069
$ws_CRP.Range("L31").Value = 12501455
$ws_CRP.Range("M31").Value = -404.2069
$ws_CRP.Range("N31").Value = -12502045

# Sheet CRP, row 34
$ws_CRP.Range("H34").Value = 7247514.5
$ws_CRP.Range("I34").Value = 699.2069
$ws_CRP.Range("J34").Value = 12501455
$ws_CRP.Range("K34").Value = 699.2069
$ws_CRP.Range("L34").Value = 12501455
$ws_CRP.Range("M34").Value = -497.2069
$ws_CRP.Range("N34").Value = -12501859

# Sheet CRP, row 41
$ws_CRP.Range("H41").Value = 9830.833000000001
$ws_CRP.Range("I41").Value = 7500
$ws_CRP.Range("J41").Value = 10996.25
$ws_CRP.Range("K41").Value = 7500
$ws_CRP.Range("L41").Value = 10996.25
$ws_CRP.Range("M41").Value = -7072
$ws_CRP.Range("N41").Value = -11852.25

# Sheet CRP, row 45
$ws_CRP.Range("H45").Value = 0
$ws_CRP.Range("I45").Value = 0
$ws_CRP.Range("K45").Value = 0
$ws_CRP.Range("M45").ClearContents()

# Sheet CRP, row 47
$ws_CRP.Range("H47").Value = 23880.666
$ws_CRP.Range("J47").Value = 23880.666
$ws_CRP.Range("L47").Value = 23880.666
$ws_CRP.Range("N47").Value = -25012.666

# Sheet CRP, row 100
$ws_CRP.Range("H100").Value = 31011.2
$ws_CRP.Range("J100").Value = 31011.2
$ws_CRP.Range("L100").Value = 31011.2
$ws_CRP.Range("N100").Value = -33175.2

$ws_GSM = $wb.Worksheets.Item("GSM")
# Sheet GSM, row 20
$ws_GSM.Range("H20").Value = 9250
$ws_GSM.Range("J20").Value = 9250
$ws_GSM.Range("L20").Value = 9250
$ws_GSM.Range("N20").Value = -9740

# Sheet GSM, row 31
$ws_GSM.Range("H31").Value = 3296
$ws_GSM.Range("I31").Value = 1620
$ws_GSM.Range("K31").Value = 1620
$ws_GSM.Range("M31").Value = -1328

# Sheet GSM, row 37
$ws_GSM.Range("H37").Value = 3296
$ws_GSM.Range("I37").Value = 1620
$ws_GSM.Range("K37").Value = 1620
$ws_GSM.Range("M37").Value = -1343

$ws_LTW = $wb.Worksheets.Item("LTW")
# Sheet LTW, row 7
$ws_LTW.Range("H7").Value = 3990
$ws_LTW.Range("I7").Value = 0
$ws_LTW.Range("K7").Value = 0
$ws_LTW.Range("M7").ClearContents()

# Sheet LTW, row 9
$ws_LTW.Range("H9").Value = 9858.5
$ws_LTW.Range("I9").Value = 175
$ws_LTW.Range("J9").Value = 14700.25
$ws_LTW.Range("K9").Value = 175
$ws_LTW.Range("L9").Value = 14700.25
$ws_LTW.Range("M9").Value = 49
$ws_LTW.Range("N9").Value = -15148.25

# Sheet LTW, row 22
$ws_LTW.Range("H22").Value = 722.08
$ws_LTW.Range("I22").Value = 400
$ws_LTW.Range("J22").Value = 766
$ws_LTW.Range("K22").Value = 400
$ws_LTW.Range("L22").Value = 766
$ws_LTW.Range("M22").Value = -105
$ws_LTW.Range("N22").Value = -1356

# Sheet LTW, row 27
$ws_LTW.Range("H27").Value = 722.08
$ws_LTW.Range("I27").Value = 400
$ws_LTW.Range("J27").Value = 766
$ws_LTW.Range("K27").Value = 400
$ws_LTW.Range("L27").Value = 766
$ws_LTW.Range("M27").Value = -293
$ws_LTW.Range("N27").Value = -980

# Sheet LTW, row 40
$ws_LTW.Range("H40").Value = 7375
$ws_LTW.Range("I40").Value = 6500
$ws_LTW.Range("J40").Value = 10000
$ws_LTW.Range("K40").Value = 6500
$ws_LTW.Range("L40").Value = 10000
$ws_LTW.Range("M40").Value = -6364
$ws_LTW.Range("N40").Value = -10272

# Sheet LTW, row 122
$ws_LTW.Range("H122").Value = 3982.5
$ws_LTW.Range("I122").Value = 3837.1428
$ws_LTW.Range("K122").Value = 11511.4284
$ws_LTW.Range("M122").Value = -9061.428400000001

# Sheet LTW, row 126
$ws_LTW.Range("H126").Value = 3990
$ws_LTW.Range("I126").Value = 0
$ws_LTW.Range("K126").Value = 0
$ws_LTW.Range("M126").ClearContents()

# Sheet LTW, row 136
$ws_LTW.Range("H136").Value = 1505.0172
$ws_LTW.Range("I136").Value = 1435.9231
$ws_LTW.Range("J136").Value = 1646.8422
$ws_LTW.Range("K136").Value = 4307.7693
$ws_LTW.Range("L136").Value = 4940.5266
$ws_LTW.Range("M136").Value = -1757.7693
$ws_LTW.Range("N136").Value = -10040.5266

$ws_WVR = $wb.Worksheets.Item("WVR")
# Sheet WVR, row 9
$ws_WVR.Range("H9").Value = 0
$ws_WVR.Range("I9").Value = 0
$ws_WVR.Range("K9").Value = 0
$ws_WVR.Range("M9").ClearContents()

# Sheet WVR, row 22
$ws_WVR.Range("H22").Value = 10000
$ws_WVR.Range("I22").Value = 0
$ws_WVR.Range("J22").Value = 10000
$ws_WVR.Range("K22").Value = 0
$ws_WVR.Range("L22").Value = 10000
$ws_WVR.Range("M22").ClearContents()
$ws_WVR.Range("N22").Value = -10586

# Sheet WVR, row 23
$ws_WVR.Range("H23").Value = 1805
$ws_WVR.Range("I23").Value = 1805
$ws_WVR.Range("J23").Value = 0
$ws_WVR.Range("K23").Value = 1805
$ws_WVR.Range("L23").Value = 0
$ws_WVR.Range("M23").Value = -1576
$ws_WVR.Range("N23").ClearContents()

# Sheet WVR, row 126
$ws_WVR.Range("H126").Value = 66668210
$ws_WVR.Range("I126").Value = 1035.7778
$ws_WVR.Range("K126").Value = 3107.3334
$ws_WVR.Range("M126").Value = -637.3334000000004
